$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "Yes"

$ws.Range("B3").Select()
